$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 675
$ws.Range("I96").Value = 675
$ws.Range("K96").Value = 2025
$ws.Range("M96").Value = -652
$ws.Range("H112").Value = 23257386
$ws.Range("J112").Value = 1736.8422
$ws.Range("L112").Value = 5210.5266
$ws.Range("N112").Value = -7426.5266
$ws.Range("H129").Value = 891.5441
$ws.Range("I129").Value = 658.3889
$ws.Range("J129").Value = 975.48
$ws.Range("K129").Value = 1975.1667
$ws.Range("L129").Value = 2926.44
$ws.Range("M129").Value = 3024.8333
$ws.Range("N129").Value = -12926.44
$ws.Range("H137").Value = 2327542.2
$ws.Range("I137").Value = 2779403.2
$ws.Range("J137").Value = 3685.5715
$ws.Range("K137").Value = 8338209.600000001
$ws.Range("L137").Value = 11056.7145
$ws.Range("M137").Value = -8335659.600000001
$ws.Range("N137").Value = -16156.7145
$ws.Range("H138").Value = 3324665
$ws.Range("I138").Value = 1117.0952
$ws.Range("J138").Value = 4875654
$ws.Range("K138").Value = 3351.2856
$ws.Range("L138").Value = 14626962
$ws.Range("M138").Value = 1788.7144
$ws.Range("N138").Value = -14637242

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24715.188
$ws.Range("I32").Value = 26766.404
$ws.Range("J32").Value = 10356.667
$ws.Range("K32").Value = 26766.404
$ws.Range("L32").Value = 10356.667
$ws.Range("M32").Value = -26479.404
$ws.Range("N32").Value = -10930.667
$ws.Range("H132").Value = 76496.92999999999
$ws.Range("I132").Value = 46859.816
$ws.Range("J132").Value = 185166.33
$ws.Range("K132").Value = 140579.448
$ws.Range("L132").Value = 555498.99
$ws.Range("M132").Value = -138049.448
$ws.Range("N132").Value = -560558.99

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1585.76
$ws.Range("I134").Value = 1443.5
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4330.5
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1795.5
$ws.Range("N134").Value = -20070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1795.186
$ws.Range("I31").Value = 1344.8572
$ws.Range("J31").Value = 2635.8
$ws.Range("K31").Value = 1344.8572
$ws.Range("L31").Value = 2635.8
$ws.Range("M31").Value = -1049.8572
$ws.Range("N31").Value = -3225.8
$ws.Range("H34").Value = 1795.186
$ws.Range("I34").Value = 1344.8572
$ws.Range("J34").Value = 2635.8
$ws.Range("K34").Value = 1344.8572
$ws.Range("L34").Value = 2635.8
$ws.Range("M34").Value = -1142.8572
$ws.Range("N34").Value = -3039.8
$ws.Range("H122").Value = 2186.875
$ws.Range("I122").Value = 1723.75
$ws.Range("J122").Value = 3113.125
$ws.Range("K122").Value = 5171.25
$ws.Range("L122").Value = 9339.375
$ws.Range("M122").Value = -2721.25
$ws.Range("N122").Value = -14239.375

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3636.3333
$ws.Range("J104").Value = 4014.5
$ws.Range("L104").Value = 12043.5
$ws.Range("N104").Value = -17285.5
$ws.Range("H113").Value = 568.3774
$ws.Range("I113").Value = 544.2857
$ws.Range("J113").Value = 572.04346
$ws.Range("K113").Value = 1632.8571
$ws.Range("L113").Value = 1716.13038
$ws.Range("M113").Value = 537.1428999999998
$ws.Range("N113").Value = -6056.130380000001
$ws.Range("H116").Value = 149140.44
$ws.Range("I116").Value = 133920
$ws.Range("J116").Value = 168166
$ws.Range("K116").Value = 401760
$ws.Range("L116").Value = 504498
$ws.Range("M116").Value = -398318
$ws.Range("N116").Value = -511382
$ws.Range("H131").Value = 879.6111
$ws.Range("J131").Value = 1032.3334
$ws.Range("L131").Value = 3097.0002
$ws.Range("N131").Value = -13177.0002
$ws.Range("H134").Value = 3317.1304
$ws.Range("I134").Value = 1190.6364
$ws.Range("J134").Value = 5266.4165
$ws.Range("K134").Value = 3571.9092
$ws.Range("L134").Value = 15799.2495
$ws.Range("M134").Value = 1498.0908
$ws.Range("N134").Value = -25939.2495

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1304.75
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H37").Value = 1304.75
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H123").Value = 45312.184
$ws.Range("J123").Value = 45312.184
$ws.Range("L123").Value = 45312.184
$ws.Range("N123").Value = -50212.184
$ws.Range("H132").Value = 155536.16
$ws.Range("I132").Value = 126005.5
$ws.Range("J132").Value = 202785.2
$ws.Range("K132").Value = 378016.5
$ws.Range("L132").Value = 608355.6000000001
$ws.Range("M132").Value = -375486.5
$ws.Range("N132").Value = -613415.6000000001
$ws.Range("N31").ClearContents()
$ws.Range("N37").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 732.25
$ws.Range("I46").Value = 643.1667
$ws.Range("K46").Value = 643.1667
$ws.Range("M46").Value = -455.1667
$ws.Range("H55").Value = 250.69048
$ws.Range("I55").Value = 221.12903
$ws.Range("J55").Value = 334
$ws.Range("K55").Value = 221.12903
$ws.Range("L55").Value = 334
$ws.Range("M55").Value = -48.12903
$ws.Range("N55").Value = -680
$ws.Range("H132").Value = 59755.11
$ws.Range("I132").Value = 3845.7778
$ws.Range("J132").Value = 115664.445
$ws.Range("K132").Value = 11537.3334
$ws.Range("L132").Value = 346993.335
$ws.Range("M132").Value = -9007.3334
$ws.Range("N132").Value = -352053.335
$ws.Range("H136").Value = 182073.25
$ws.Range("I136").Value = 85681.836
$ws.Range("J136").Value = 471247.5
$ws.Range("K136").Value = 257045.508
$ws.Range("L136").Value = 1413742.5
$ws.Range("M136").Value = -254495.508
$ws.Range("N136").Value = -2825100

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1784.2368
$ws.Range("I122").Value = 1547.7778
$ws.Range("J122").Value = 2364.6365
$ws.Range("K122").Value = 4643.3334
$ws.Range("L122").Value = 7093.9095
$ws.Range("M122").Value = -2193.3334
$ws.Range("N122").Value = -11993.9095
$ws.Range("H126").Value = 1358.2759
$ws.Range("I126").Value = 942.8261
$ws.Range("J126").Value = 2950.8333
$ws.Range("K126").Value = 2828.4783
$ws.Range("L126").Value = 8852.499899999999
$ws.Range("M126").Value = -358.4782999999998
$ws.Range("N126").Value = -13792.4999
$ws.Range("H132").Value = 63742.97
$ws.Range("I132").Value = 48414.094
$ws.Range("J132").Value = 93007.17999999999
$ws.Range("K132").Value = 145242.282
$ws.Range("L132").Value = 279021.54
$ws.Range("M132").Value = -142712.282
$ws.Range("N132").Value = -284081.54
$ws.Range("H136").Value = 334083.34
$ws.Range("I136").Value = 250875
$ws.Range("J136").Value = 500500
$ws.Range("K136").Value = 752625
$ws.Range("L136").Value = 1501500
$ws.Range("M136").Value = -750075
$ws.Range("N136").Value = -1506600
